$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are numeric-looking strings need to be forced to
# Text format first, otherwise Excel auto-converts them to numbers and
# precision/trailing-zero formatting (e.g. "1.110", "1.000") would be lost.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.185.30"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").Value = "1.827.57"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("D5").Value = "311.02"
$ws.Range("E5").Value = "  -0.80%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "0.4959"
$ws.Range("E7").Value = "  -3.58%  "
$ws.Range("D8").Value = "0.3932"
$ws.Range("E8").Value = "  -1.54%  "
$ws.Range("D9").Value = "0.09898"
$ws.Range("E9").Value = "  +25.64%  "
$ws.Range("D10").Value = "1.110"
$ws.Range("D11").Value = "41.19"
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("D12").Value = "6.447"
$ws.Range("E12").Value = "  +0.94%  "
$ws.Range("D13").Value = "20.64"
$ws.Range("E13").Value = "  +1.08%  "
$ws.Range("D14").Value = "1.001"
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("D15").Value = "1.829.06"
$ws.Range("E15").Value = "  +1.39%  "
$ws.Range("D16").Value = "7.313"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").Value = "0.00001145"
$ws.Range("E17").Value = "  +5.57%  "
$ws.Range("D18").Value = "92.91"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").Value = "0.06652"
$ws.Range("E19").Value = "  +1.14%  "
$ws.Range("D21").Value = "17.25"
$ws.Range("E21").Value = "  -0.67%  "
$ws.Range("D22").Value = "6.003"
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("D23").Value = "28.229.96"
$ws.Range("E23").Value = "  -0.59%  "
$ws.Range("D24").Value = "11.35"
$ws.Range("E24").Value = "  +1.32%  "
$ws.Range("D25").Value = "2.243"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").Value = "158.87"
$ws.Range("E26").Value = "  -1.40%  "
$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").Value = "2.036.13"
$ws.Range("E27").Value = "  +0.88%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "20.75"
$ws.Range("E28").Value = "  +1.00%  "
$ws.Range("D29").Value = "2.425"
$ws.Range("E29").Value = "  +0.89%  "
$ws.Range("E30").Value = "  -1.36%  "
$ws.Range("E31").Value = "  -3.06%  "
$ws.Range("D32").Value = "1.042"
$ws.Range("E32").Value = "  -2.62%  "
$ws.Range("D33").Value = "5.601"
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("D34").Value = "3.607"
$ws.Range("E34").Value = "  -1.52%  "
$ws.Range("D35").Value = "0.06756"
$ws.Range("E35").Value = "  -6.94%  "
$ws.Range("D36").Value = "9.059"
$ws.Range("E36").Value = "  -1.07%  "
$ws.Range("D37").Value = "0.02345"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").Value = "0.2155"
$ws.Range("E38").Value = "  -1.21%  "
$ws.Range("D39").Value = "4.983"
$ws.Range("E39").Value = "  -1.99%  "
$ws.Range("D40").Value = "11.40"
$ws.Range("E40").Value = "  -2.22%  "
$ws.Range("D41").Value = "0.6225"
$ws.Range("E41").Value = "  +0.21%  "
$ws.Range("E42").Value = "  +1.53%  "
$ws.Range("D43").Value = "1.000"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").Value = "13.18"
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("D45").Value = "0.5943"
$ws.Range("E45").Value = "  -1.11%  "
$ws.Range("D46").Value = "3.702"
$ws.Range("E46").Value = "  -1.10%  "
$ws.Range("D47").Value = "1.271"
$ws.Range("E47").Value = "  -3.10%  "
$ws.Range("D48").Value = "124.26"
$ws.Range("E48").Value = "  -1.17%  "
$ws.Range("D49").Value = "1.948"
$ws.Range("E49").Value = "  +0.58%  "
$ws.Range("D50").Value = "1.183"
$ws.Range("D51").Value = "0.06794"
$ws.Range("E51").Value = "  -0.85%  "
